$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 to row 46 (row 1 is the header: entity_id/type/value).
$lastRow = 46

for ($r = 2; $r -le $lastRow; $r++) {
    $typeCell = $ws.Cells.Item($r, 2)
    $valueCell = $ws.Cells.Item($r, 3)

    $oldType = $typeCell.Value()
    $oldValue = $valueCell.Value()

    if ([string]::IsNullOrEmpty($oldType) -or [string]::IsNullOrEmpty($oldValue)) {
        continue
    }

    # association: club-sports -> club, uil-sports -> uil
    if ($oldType -eq "club-sports") {
        $assoc = "club"
    } elseif ($oldType -eq "uil-sports") {
        $assoc = "uil"
    } else {
        continue
    }

    # gender suffix comes from the old value "Softball-<Gender>"
    $parts = $oldValue -split "-"
    if ($parts.Length -lt 2) {
        continue
    }
    $gender = $parts[1].ToLower()

    $newType = "sports_" + $assoc + "_" + $gender
    $newValue = "Softball"

    $typeCell.Value = $newType
    $valueCell.Value = $newValue
}
